$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue 'D2' '64.052.01'
Set-TextValue 'E2' '  +4.16%  '

# Row 3
Set-TextValue 'D3' '2.786.22'
Set-TextValue 'E3' '  +4.90%  '

# Row 4
Set-TextValue 'E4' '  +0.00%  '

# Row 5
Set-TextValue 'D5' '584.46'
Set-TextValue 'E5' '  +0.52%  '

# Row 6
Set-TextValue 'D6' '161.52'
Set-TextValue 'E6' '  +11.57%  '

# Row 7
Set-TextValue 'D7' '0.625'
Set-TextValue 'E7' '  +3.67%  '

# Row 8
Set-TextValue 'D8' '0.995'

# Row 9
Set-TextValue 'D9' '2.804.55'
Set-TextValue 'E9' '  +4.95%  '

# Row 10
Set-TextValue 'D10' '6.87'
Set-TextValue 'E10' '  +3.90%  '

# Row 11
Set-TextValue 'D11' '0.114'
Set-TextValue 'E11' '  +3.84%  '

# Row 12
Set-TextValue 'D12' '0.399'
Set-TextValue 'E12' '  +4.21%  '

# Row 13
Set-TextValue 'E13' '  +1.01%  '

# Row 14
Set-TextValue 'D14' '3.259.76'
Set-TextValue 'E14' '  +4.36%  '

# Row 15
Set-TextValue 'D15' '27.91'
Set-TextValue 'E15' '  +6.87%  '

# Row 16
Set-TextValue 'D16' '64.006.19'
Set-TextValue 'E16' '  +4.29%  '

# Row 17
Set-TextValue 'D17' '0.0000160'
Set-TextValue 'E17' '  +8.61%  '

# Row 18
Set-TextValue 'D18' '2.779.54'
Set-TextValue 'E18' '  +4.45%  '

# Row 19
Set-TextValue 'D19' '12.35'
Set-TextValue 'E19' '  +5.70%  '

# Row 20
Set-TextValue 'D20' '5.01'
Set-TextValue 'E20' '  +4.87%  '

# Row 21
Set-TextValue 'D21' '368.16'
Set-TextValue 'E21' '  +3.55%  '

# Row 22
Set-TextValue 'D22' '7.10'
Set-TextValue 'E22' '  +2.65%  '

# Row 23
Set-TextValue 'D23' '0.548'
Set-TextValue 'E23' '  +4.31%  '

# Row 24
Set-TextValue 'D24' '1.01'
Set-TextValue 'E24' '  +0.70%  '

# Row 25
Set-TextValue 'D25' '67.63'
Set-TextValue 'E25' '  +4.63%  '

# Row 26
Set-TextValue 'D26' '0.174'
Set-TextValue 'E26' '  +5.84%  '

# Row 27
Set-TextValue 'D27' '8.70'
Set-TextValue 'E27' '  +2.44%  '

# Row 28
Set-TextValue 'B28' 'PEPE'
Set-TextValue 'C28' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D28' '0.0₃0966'
Set-TextValue 'E28' '  +17.31%  '

# Row 29
Set-TextValue 'B29' 'Binance-PegBSC-USD'
Set-TextValue 'C29' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D29' '1.00'
Set-TextValue 'E29' '  +0.38%  '

# Row 30
Set-TextValue 'D30' '2.04'
Set-TextValue 'E30' '  +1.79%  '

# Row 31
Set-TextValue 'D31' '7.31'
Set-TextValue 'E31' '  +5.45%  '

# Row 32
Set-TextValue 'E32' '  +12.18%  '

# Row 33
Set-TextValue 'D33' '173.69'
Set-TextValue 'E33' '  +2.41%  '

# Row 34
Set-TextValue 'B34' 'EthereumClassic'
Set-TextValue 'C34' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D34' '20.94'
Set-TextValue 'E34' '  +3.74%  '

# Row 35
Set-TextValue 'B35' 'USDe'
Set-TextValue 'C35' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D35' '0.997'
Set-TextValue 'E35' '  -0.08%  '

# Row 36
Set-TextValue 'D36' '5.05'
Set-TextValue 'E36' '  +7.92%  '

# Row 37
Set-TextValue 'D37' '1.49'
Set-TextValue 'E37' '  +7.60%  '

# Row 38
Set-TextValue 'D38' '1.85'
Set-TextValue 'E38' '  +7.06%  '

# Row 39
Set-TextValue 'D39' '1.03'
Set-TextValue 'E39' '  +2.66%  '

# Row 40
Set-TextValue 'B40' 'Bittensor'
Set-TextValue 'C40' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D40' '344.01'
Set-TextValue 'E40' '  +0.24%  '

# Row 41
Set-TextValue 'B41' 'Filecoin'
Set-TextValue 'C41' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D41' '4.29'
Set-TextValue 'E41' '  +3.19%  '

# Row 42
Set-TextValue 'D42' '6.27'
Set-TextValue 'E42' '  +15.61%  '

# Row 43
Set-TextValue 'D43' '39.76'
Set-TextValue 'E43' '  +3.28%  '

# Row 44
Set-TextValue 'D44' '22.64'
Set-TextValue 'E44' '  +9.35%  '

# Row 45
Set-TextValue 'D45' '22.69'
Set-TextValue 'E45' '  +6.98%  '

# Row 46
Set-TextValue 'D46' '0.0610'
Set-TextValue 'E46' '  +5.24%  '

# Row 47
Set-TextValue 'B47' 'Mantle'
Set-TextValue 'C47' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D47' '0.654'
Set-TextValue 'E47' '  +3.96%  '

# Row 48
Set-TextValue 'B48' 'VeChain'
Set-TextValue 'C48' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D48' '0.0261'
Set-TextValue 'E48' '  +2.72%  '

# Row 49
Set-TextValue 'D49' '138.71'
Set-TextValue 'E49' '  +2.10%  '

# Row 50
Set-TextValue 'D50' '0.103'
Set-TextValue 'E50' '  +2.84%  '

# Row 51
Set-TextValue 'D51' '2.168.57'
Set-TextValue 'E51' '  +3.13%  '
